# Add a "Length" column (F) to the "XX" worksheet, with per-variable
# field lengths, mirroring the addition of a Length column to the
# dataset specification metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XX")

$ws.Range("F1").Value = "Length"

$ws.Range("F2").Value  = 200   # STUDYID
$ws.Range("F3").Value  = 200   # DOMAIN
$ws.Range("F4").Value  = 200   # USUBJID
$ws.Range("F5").Value  = 8     # XXSEQ
$ws.Range("F6").Value  = 8     # XXTESTCD
$ws.Range("F7").Value  = 40    # XXTEST
$ws.Range("F8").Value  = 200   # XXORRES
$ws.Range("F9").Value  = 1     # XXBLFL
$ws.Range("F10").Value = 200   # VISIT
$ws.Range("F11").Value = 200   # EPOCH
$ws.Range("F12").Value = 19    # XXDTC
$ws.Range("F13").Value = 8     # XXDY
